# Auto-generated edit script: applies per-cell numeric corrections
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR 'Typhon_Profits' market-data
# sheets, per the scheduled-runner diff (H..N columns: current market
# price / NQ / HQ averages, leve prices, and leve profit figures).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 293.27274
$ws.Range("I4").Value = 155.5
$ws.Range("K4").Value = 155.5
$ws.Range("M4").Value = -41.5
$ws.Range("H100").Value = 1905.5264
$ws.Range("I100").Value = 687.375
$ws.Range("J100").Value = 2791.4546
$ws.Range("K100").Value = 687.375
$ws.Range("L100").Value = 2791.4546
$ws.Range("M100").Value = -146.375
$ws.Range("N100").Value = -3873.4546
$ws.Range("H101").Value = 3000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 3000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 9000
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -12244
$ws.Range("H103").Value = 1250250
$ws.Range("I103").Value = 2500000
$ws.Range("K103").Value = 7500000
$ws.Range("M103").Value = -7499414
$ws.Range("H116").Value = 4938.077
$ws.Range("I116").Value = 2438.8
$ws.Range("J116").Value = 6500.125
$ws.Range("K116").Value = 2438.8
$ws.Range("L116").Value = 6500.125
$ws.Range("M116").Value = 1003.2
$ws.Range("N116").Value = -13384.125
$ws.Range("H121").Value = 1979.4445
$ws.Range("J121").Value = 2345
$ws.Range("L121").Value = 7035
$ws.Range("N121").Value = -10529
$ws.Range("H129").Value = 271013.1
$ws.Range("I129").Value = 292.42856
$ws.Range("J129").Value = 334181.22
$ws.Range("K129").Value = 877.28568
$ws.Range("L129").Value = 1002543.66
$ws.Range("M129").Value = 4122.71432
$ws.Range("N129").Value = -1012543.66
$ws.Range("H132").Value = 5152.875
$ws.Range("I132").Value = 5429.7334
$ws.Range("K132").Value = 16289.2002
$ws.Range("M132").Value = -13759.2002
$ws.Range("H135").Value = 25008134
$ws.Range("I135").Value = 510.86667
$ws.Range("K135").Value = 4597.80003
$ws.Range("M135").Value = -2062.80003
$ws.Range("H137").Value = 1670.0416
$ws.Range("I137").Value = 1384.05
$ws.Range("K137").Value = 4152.15
$ws.Range("M137").Value = -1602.15
$ws.Range("H138").Value = 1579.31
$ws.Range("I138").Value = 630.1613
$ws.Range("J138").Value = 2005.7391
$ws.Range("K138").Value = 1890.4839
$ws.Range("L138").Value = 6017.2173
$ws.Range("M138").Value = 3249.5161
$ws.Range("N138").Value = -16297.2173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 142858860
$ws.Range("I74").Value = 250000900
$ws.Range("K74").Value = 250000900
$ws.Range("M74").Value = -250000026
$ws.Range("H77").Value = 142858860
$ws.Range("I77").Value = 250000900
$ws.Range("K77").Value = 1250004500
$ws.Range("M77").Value = -1250000132
$ws.Range("H122").Value = 2148.4666
$ws.Range("I122").Value = 2017.4615
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6052.3845
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3602.3845
$ws.Range("N122").Value = -13900
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4905.7827
$ws.Range("I134").Value = 5401.65
$ws.Range("K134").Value = 16204.95
$ws.Range("M134").Value = -13669.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H31").Value = 13733.034
$ws.Range("I31").Value = 21685.4
$ws.Range("J31").Value = 5212.643
$ws.Range("K31").Value = 21685.4
$ws.Range("L31").Value = 5212.643
$ws.Range("M31").Value = -21390.4
$ws.Range("N31").Value = -5802.643
$ws.Range("H34").Value = 13733.034
$ws.Range("I34").Value = 21685.4
$ws.Range("J34").Value = 5212.643
$ws.Range("K34").Value = 21685.4
$ws.Range("L34").Value = 5212.643
$ws.Range("M34").Value = -21483.4
$ws.Range("N34").Value = -5616.643
$ws.Range("H58").Value = 37388.355
$ws.Range("I58").Value = 2287.4285
$ws.Range("J58").Value = 72489.28999999999
$ws.Range("K58").Value = 2287.4285
$ws.Range("L58").Value = 72489.28999999999
$ws.Range("M58").Value = -2084.4285
$ws.Range("N58").Value = -72895.28999999999
$ws.Range("H134").Value = 1567.0454
$ws.Range("I134").Value = 1040.0714
$ws.Range("J134").Value = 2489.25
$ws.Range("K134").Value = 3120.2142
$ws.Range("L134").Value = 7467.75
$ws.Range("M134").Value = -585.2142000000003
$ws.Range("N134").Value = -12537.75
$ws.Range("H136").Value = 37388.355
$ws.Range("I136").Value = 2287.4285
$ws.Range("J136").Value = 72489.28999999999
$ws.Range("K136").Value = 6862.2855
$ws.Range("L136").Value = 217467.87
$ws.Range("M136").Value = -4312.2855
$ws.Range("N136").Value = -222567.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 526.55554
$ws.Range("I113").Value = 443.63635
$ws.Range("J113").Value = 656.8570999999999
$ws.Range("K113").Value = 1330.90905
$ws.Range("L113").Value = 1970.5713
$ws.Range("M113").Value = 839.09095
$ws.Range("N113").Value = -6310.5713
$ws.Range("H131").Value = 223019.83
$ws.Range("I131").Value = 1020
$ws.Range("J131").Value = 238876.95
$ws.Range("K131").Value = 3060
$ws.Range("L131").Value = 716630.8500000001
$ws.Range("M131").Value = 1980
$ws.Range("N131").Value = -726710.8500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3152.5
$ws.Range("I80").Value = 2790.3845
$ws.Range("J80").Value = 3429.4119
$ws.Range("K80").Value = 2790.3845
$ws.Range("L80").Value = 3429.4119
$ws.Range("M80").Value = -1792.3845
$ws.Range("N80").Value = -5425.4119
$ws.Range("H83").Value = 3152.5
$ws.Range("I83").Value = 2790.3845
$ws.Range("J83").Value = 3429.4119
$ws.Range("K83").Value = 13951.9225
$ws.Range("L83").Value = 17147.0595
$ws.Range("M83").Value = -8959.922500000001
$ws.Range("N83").Value = -27131.0595
$ws.Range("H102").Value = 3519
$ws.Range("I102").Value = 3519
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3519
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1897
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 32856.777
$ws.Range("I132").Value = 5294.933
$ws.Range("J132").Value = 170666
$ws.Range("K132").Value = 15884.799
$ws.Range("L132").Value = 511998
$ws.Range("M132").Value = -13354.799
$ws.Range("N132").Value = -517058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2286.682
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 2513.8
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2513.8
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -4011.8
$ws.Range("H71").Value = 2286.682
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 2513.8
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 12569
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -20057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3399.8
$ws.Range("J15").Value = 3399.8
$ws.Range("L15").Value = 3399.8
$ws.Range("N15").Value = -3975.8
$ws.Range("H100").Value = 160.2
$ws.Range("J100").Value = 190
$ws.Range("L100").Value = 380
$ws.Range("N100").Value = -1462
$ws.Range("H126").Value = 1142
$ws.Range("I126").Value = 904
$ws.Range("J126").Value = 1163.6364
$ws.Range("K126").Value = 2712
$ws.Range("L126").Value = 3490.9092
$ws.Range("M126").Value = -242
$ws.Range("N126").Value = -8430.9092

Write-Output "Applied $(190) value updates and $(2) clears across 8 sheets."